$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# U1: new header date "06-10-2020" matching the style of N1:T1 (bold, thin border, centered/top)
$u1 = $ws.Range("U1")
$u1.NumberFormat = "@"
$u1.Value = "06-10-2020"
$u1.Font.Bold = $true
$u1.Borders.LineStyle = 1
$u1.HorizontalAlignment = -4108
$u1.VerticalAlignment = -4160

# U2:U36: new numeric data column (active-case counts for 06-10-2020)
$ws.Range("U2").Value = 186
$ws.Range("U3").Value = 51060
$ws.Range("U4").Value = 2989
$ws.Range("U5").Value = 33467
$ws.Range("U6").Value = 11523
$ws.Range("U7").Value = 1604
$ws.Range("U8").Value = 27857
$ws.Range("U9").Value = 99
$ws.Range("U10").Value = 23080
$ws.Range("U11").Value = 4803
$ws.Range("U12").Value = 16718
$ws.Range("U13").Value = 11822
$ws.Range("U14").Value = 3156
$ws.Range("U15").Value = 14696
$ws.Range("U16").Value = 10436
$ws.Range("U17").Value = 115496
$ws.Range("U18").Value = 84958
$ws.Range("U19").Value = 1166
$ws.Range("U20").Value = 18757
$ws.Range("U21").Value = 252721
$ws.Range("U22").Value = 2696
$ws.Range("U23").Value = 2217
$ws.Range("U24").Value = 291
$ws.Range("U25").Value = 1155
$ws.Range("U26").Value = 28006
$ws.Range("U27").Value = 4513
$ws.Range("U28").Value = 12895
$ws.Range("U29").Value = 21215
$ws.Range("U30").Value = 598
$ws.Range("U31").Value = 45881
$ws.Range("U32").Value = 26644
$ws.Range("U33").Value = 4876
$ws.Range("U34").Value = 8701
$ws.Range("U35").Value = 45024
$ws.Range("U36").Value = 27717
